$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cargo header cells in row 1 (H1, I1)
$ws.Range("H1").Value = "Груз 1"
$ws.Range("I1").Value = "Груз 2"

# Update driver row 3 with corrected data
$ws.Range("C3").Value = "Илья"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "19"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "Вольво"
$ws.Range("F3").Value = "Вроцлав"

$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "20"
$ws.Range("G3").Style = "Normal"

# Assign new cargo columns to driver in row 3
$ws.Range("J3").Value = "Москва_Ростов-на-Дону"
$ws.Range("K3").Value = "Керчь_Новочеркасск"

# Update selection/view to reflect newly added cargo columns
$ws.Range("J3:M3").Select()
$excel.ActiveWindow.ScrollColumn = 6
